$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsCodebook = $wb.Worksheets.Item("Codebook")

# --- Data sheet: add new columns D (hours_sleep) and E (eye_color) ---
$wsData.Cells.Item(1,4).Value = "hours_sleep"
$wsData.Cells.Item(1,5).Value = "eye_color"

$wsData.Cells.Item(2,4).Value = 6
$wsData.Cells.Item(2,5).Value = "brown "

$wsData.Cells.Item(3,4).Value = 7
$wsData.Cells.Item(3,5).Value = "blue"

$wsData.Cells.Item(4,4).Value = 7.5
$wsData.Cells.Item(4,5).Value = "green"

$wsData.Cells.Item(5,4).Value = 8
$wsData.Cells.Item(5,5).Value = "brown "

$wsData.Cells.Item(6,4).Value = 9
$wsData.Cells.Item(6,5).Value = "blue"

$wsData.Cells.Item(7,4).Value = 5
$wsData.Cells.Item(7,5).Value = "blue"

$wsData.Cells.Item(8,4).Value = 8
$wsData.Cells.Item(8,5).Value = "brown "

$wsData.Cells.Item(9,4).Value = 9
$wsData.Cells.Item(9,5).Value = "green"

$wsData.Cells.Item(10,4).Value = 7
$wsData.Cells.Item(10,5).Value = "green"

$wsData.Cells.Item(11,4).Value = 7
$wsData.Cells.Item(11,5).Value = "brown "

$wsData.Cells.Item(12,4).Value = 6.5
$wsData.Cells.Item(12,5).Value = "brown "

$wsData.Cells.Item(13,4).Value = 6
$wsData.Cells.Item(13,5).Value = "brown "

$wsData.Cells.Item(14,4).Value = 5
$wsData.Cells.Item(14,5).Value = "blue"

$wsData.Cells.Item(15,4).Value = 6
$wsData.Cells.Item(15,5).Value = "blue"

# --- Codebook sheet: document the two new variables ---
$wsCodebook.Cells.Item(5,1).Value = "hours_sleep"
$wsCodebook.Cells.Item(5,2).Value = "numbers of hours a person sleeps per night"
$wsCodebook.Cells.Item(5,3).Value = "numeric (hours)"

$wsCodebook.Cells.Item(6,1).Value = "eye_color"
$wsCodebook.Cells.Item(6,2).Value = "eye color category"
$wsCodebook.Cells.Item(6,3).Value = "brown, blue, green"

# --- View state: selection moves to one past the last data row on Data,
#     then the Codebook sheet becomes the active tab with its new last cell selected ---
$wsData.Activate()
$wsData.Range("E16").Select()

$wsCodebook.Activate()
$wsCodebook.Range("C6").Select()
